# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# --- Insert a new sheet "2022-Q1" right before the "总计" (totals) sheet ---
$totalSheet = $wb.Worksheets.Item(4)
$q1Sheet = $wb.Worksheets.Add($totalSheet)
$q1Sheet.Name = "2022-Q1"

# Reuse the header/row styling already present on the "2021-Q4" sheet so the
# new sheet's header row + index column formatting matches the rest of the
# workbook (bold, centered, thin border).
$styleSrc = $wb.Worksheets.Item(3)
$styleSrc.Range("A1:H2").Copy($q1Sheet.Range("A1"))
$styleSrc.Range("A2:H2").Copy($q1Sheet.Range("A3"))
# A1 is empty in the source range too, but the copy still materialises a
# (valueless) cell record for it; drop that so the sheet matches the usual
# "row 1 only has B1:H1" shape used elsewhere in this workbook.
$q1Sheet.Range("A1").ClearContents()

$q1Sheet.Range("D1").Value = "基金规模"

$q1Sheet.Range("A2").Value = 0
$q1Sheet.Range("B2").Value = "'004685"
$q1Sheet.Range("C2").Value = "金元顺安元启灵活配置混合"
$q1Sheet.Range("D2").Value = "'5.00"
$q1Sheet.Range("E2").Value = "'75.79"
$q1Sheet.Range("F2").Value = "'0.98"
$q1Sheet.Range("G2").Value = "'0.0490"
$q1Sheet.Range("H2").Value = 9
# Drop the quote-prefix formatting picked up from the text assignments above
# so these cells stay on the plain/default style like the source sheet.
$q1Sheet.Range("B2:G2").ClearFormats()

$q1Sheet.Range("A3").Value = 1
$q1Sheet.Range("B3").Value = "'005126"
$q1Sheet.Range("C3").Value = "银河量化稳进混合"
$q1Sheet.Range("D3").Value = "'0.10"
$q1Sheet.Range("E3").Value = "'78.20"
$q1Sheet.Range("F3").Value = "'2.18"
$q1Sheet.Range("G3").Value = "'0.0022"
$q1Sheet.Range("H3").Value = 4
$q1Sheet.Range("B3:G3").ClearFormats()

# --- Add the 2022-Q1 summary row to the top of the "总计" sheet ---
$totalSheet = $wb.Worksheets.Item(5)
$totalSheet.Rows(2).Insert()

# New blank row doesn't inherit the index column's style; copy it from the
# row right below (which used to be row 2, still styled).
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.05
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
